$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / Row 7: new "Electrolyte Step Size" / "Number of Electrolyte " labels ---
# A6 previously held "=D3/6"; clear its old number style back to Normal before
# writing the new label so it ends up with the default (no explicit) style.
$ws.Range("A6").Style = "Normal"
$ws.Range("A6").Value = "Electrolyte Step Size"
$ws.Range("B6").Formula = "=D3/B7"

$ws.Range("A7").Value = "Number of Electrolyte "
$ws.Range("B7").Value = 5

# --- Row 10: step sequence now derived from $B$6 (electrolyte step size) ---
$ws.Range("A10:H10").NumberFormat = "0.000000E+00"

$ws.Range("A10").Value = 0
$ws.Range("B10").Formula = "=A10+`$B`$6"
$ws.Range("C10").Formula = "=B10+`$B`$6"
$ws.Range("D10").Formula = "=C10+`$B`$6"
$ws.Range("E10").Formula = "=D10+`$B`$6"
$ws.Range("F10").Formula = "=E10+`$B`$6"
$ws.Range("G10").Formula = "=F10+`$B`$6"
$ws.Range("H10").Formula = "=G10+`$B`$6"

# New helper cell I10 (kept empty, same numeric style used elsewhere on the sheet)
$ws.Range("I10").NumberFormat = "0.00E+00"

# --- Column widths: A wider for the new labels, B:H sized to the step values ---
$ws.Columns("A").ColumnWidth = 17.42
$ws.Range("B:H").ColumnWidth = 11.5

# --- Selection / view ---
[void]$ws.Range("E8").Select()
